$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Green color used for the "Mini Project Completed" row (RGB 00B050 -> OLE 5287936)
$green = 5287936

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $c1 = $t.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)
    $c2 = $t.Cell($r, 2).Range.Text.TrimEnd([char]13, [char]7)
    $c3 = $t.Cell($r, 3).Range.Text.TrimEnd([char]13, [char]7)

    if ($c1 -eq "17" -and $c2 -eq "POC" -and $c3 -eq "Proof Of Concept") {
        $t.Rows.Item($r).Range.Font.Color = $green
    }
}
